$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "19. Aug: " paragraph (Paragraph.Range.Text includes the
# trailing paragraph-mark char, so compare against the text without it).
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text
    $noMark = $full.Substring(0, $full.Length - 1)
    if ($noMark -eq "19. Aug: ") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# Insertion point = right before the paragraph mark (i.e. right after
# the existing "19. Aug: " run), so new text lands inside the same
# paragraph, before the bookmark.
$insStart = $p.Range.End - 1

$ins = $d.Range($insStart, $insStart)
$ins.InsertAfter("4:30-5:15")

# Force the newly inserted text to live in its own runs (matching the
# target markup) by toggling a boolean character property on the
# relevant sub-ranges; Bold defaults back to "not set" so this leaves
# no residual formatting once reset to 0.
$rAll = $d.Range($insStart, $insStart + 9)
$rAll.Bold = 1
$rAll.Bold = 0

$rTime2 = $d.Range($insStart + 5, $insStart + 9)
$rTime2.Bold = 1
$rTime2.Bold = 0

# ------------------------------------------------------------------
# Insert the three new paragraphs right after the "19. Aug: " paragraph
# ------------------------------------------------------------------

# Paragraph 1: camera-bug fix note (single run)
$p.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($targetIndex + 1)
$p1.Range.Text = "- fixed camera bug: camera was following with lag. Had to call Update() instead of FixedUpdate for moving the camera."

# Paragraph 2: movement-bug fix note (two runs: "- " + rest)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($targetIndex + 2)
$p2Start = $p2.Range.Start
$p2.Range.Text = "- fixed movement bug: movement works now in multiplayer"

$rDash = $d.Range($p2Start, $p2Start + 2)
$rDash.Bold = 1
$rDash.Bold = 0

# Paragraph 3: TODO note followed by a manual line break in its own run
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($targetIndex + 3)
$p3Start = $p3.Range.Start
$todoText = "- TODO: deactivate velocity force on collision!!"
$p3.Range.Text = $todoText

$brPos = $p3Start + $todoText.Length
$rBrIns = $d.Range($brPos, $brPos)
$rBrIns.InsertAfter([char]11)

$rBr = $d.Range($brPos, $brPos + 1)
$rBr.Bold = 1
$rBr.Bold = 0
